$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three cells whose "- modif HHhMM" suffix was removed,
#     leaving a trailing space behind (per the shared-string diff). ---
$ws.Range("A1").Value = "Donnée A1 "
$ws.Range("B2").Value = "Donnée B2 "
$ws.Range("C3").Value = "Donnée C3 "

# --- Insert a Line sparkline in C11 sourced from column A. ---
$ws.Range("C11").SparklineGroups.Add(7, "Feuil1!A:A") | Out-Null

# --- Move the active selection from C3 to F5. ---
$ws.Range("F5").Select() | Out-Null
